# [MIG] 12.0 account_bank_statement_import_adyen, account_bank_statement_clearing_account
#
# Update the Adyen bank-statement-import test fixture:
#  - the sample "Gross/Net Currency" code used throughout the sheet changes
#    from EUR to USD
#  - the Gross Credit (GC) amount on the "Settled" sample row (M10) is
#    updated from 666 to 1598
#  - refresh the active selection on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell that shows the currency code "EUR" (Gross Currency / Net
# Currency columns) becomes "USD".
$ws.Cells.Replace("EUR", "USD") | Out-Null

# Gross Credit (GC) amount for the Settled/mc sample row.
$ws.Range("M10").Value = 1598

# Leave the selection the way the refreshed sheet had it.
$ws.Range("L9").Select() | Out-Null
